$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.039.93'
$ws.Range("E2").Value = '  +0.47%  '

$ws.Range("E3").Value = '  +0.62%  '

$ws.Range("E4").Value = '  +0.31%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.45'
$ws.Range("E5").Value = '  +0.94%  '

$ws.Range("E6").Value = '  +0.63%  '

$ws.Range("E7").Value = '  +0.29%  '

$ws.Range("E9").Value = '  +1.31%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.66'
$ws.Range("E10").Value = '  +0.57%  '

$ws.Range("E11").Value = '  +0.58%  '

$ws.Range("D12").Value = '1.698.63'
$ws.Range("E12").Value = '  +4.32%  '

$ws.Range("D13").Value = '1.872.28'
$ws.Range("E13").Value = '  +0.65%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.28'
$ws.Range("E14").Value = '  +0.98%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.544'
$ws.Range("E15").Value = '  +0.37%  '

$ws.Range("E16").Value = '  +1.32%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.35'
$ws.Range("E17").Value = '  +0.94%  '

$ws.Range("D18").Value = '26.051.79'
$ws.Range("E18").Value = '  +0.54%  '

$ws.Range("E19").Value = '  +0.36%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '195.54'
$ws.Range("E20").Value = '  +1.43%  '

$ws.Range("E21").Value = '  -0.54%  '

$ws.Range("E22").Value = '  +0.04%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.24'
$ws.Range("E23").Value = '  +0.15%  '

$ws.Range("B24").Value = 'Stellar'
$ws.Range("C24").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.132'
$ws.Range("E24").Value = '  +5.19%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.80'
$ws.Range("E25").Value = '  +0.09%  '

$ws.Range("B26").Value = 'BinanceUSD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.01'
$ws.Range("E26").Value = '  +0.71%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '143.80'
$ws.Range("E27").Value = '  +0.45%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.90'
$ws.Range("E28").Value = '  +0.73%  '

$ws.Range("E29").Value = '  +0.68%  '

$ws.Range("E30").Value = '  +1.21%  '

$ws.Range("E31").Value = '  +0.05%  '

$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.29'
$ws.Range("E32").Value = '  -0.14%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.27'
$ws.Range("E33").Value = '  +1.36%  '

$ws.Range("E34").Value = '  -2.68%  '

$ws.Range("E35").Value = '  +1.19%  '

$ws.Range("E36").Value = '  +0.62%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '431.42'
$ws.Range("E37").Value = '  +20.48%  '

$ws.Range("D38").Value = '1.134.23'
$ws.Range("E38").Value = '  -0.31%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.541'
$ws.Range("E39").Value = '  -1.42%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.46'
$ws.Range("E40").Value = '  -0.49%  '

$ws.Range("E41").Value = '  +0.23%  '

$ws.Range("E42").Value = '  +0.97%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.39'
$ws.Range("E43").Value = '  +0.16%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.799'
$ws.Range("E44").Value = '  -0.86%  '

$ws.Range("D45").Value = '1.781.55'
$ws.Range("E45").Value = '  +0.68%  '

$ws.Range("E46").Value = '  +3.72%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.78'
$ws.Range("E47").Value = '  +0.93%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0525'
$ws.Range("E48").Value = '  +0.42%  '

$ws.Range("E49").Value = '  +0.65%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.80'
$ws.Range("E50").Value = '  +2.35%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.414'
$ws.Range("E51").Value = '  -0.19%  '
